$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.543891666666667
$ws.Range("H2").Value = 7.631675
$ws.Range("I2").Value = 0.0487891203504456
$ws.Range("J2").Value = 0.04878912035044559
$ws.Range("M2").Value = 40.81054266666667
$ws.Range("N2").Value = 122.431628
$ws.Range("O2").Value = 0.9943414173631485
$ws.Range("P2").Value = 0.9943414173631485
$ws.Range("Q2").Value = 103.8175994018778
$ws.Range("R2").Value = 934.3583946169
$ws.Range("S2").Value = 0.04851304308116331
$ws.Range("T2").Value = 0.0485130430811633

# Row 3
$ws.Range("G3").Value = 2.543891666666667
$ws.Range("H3").Value = 7.631675
$ws.Range("I3").Value = 0.0487891203504456
$ws.Range("J3").Value = 0.04878912035044559
$ws.Range("O3").Value = 0.0002749163555820933
$ws.Range("P3").Value = 0.0002749163555820933
$ws.Range("Q3").Value = 0.02870357763888889
$ws.Range("R3").Value = 0.25833219875
$ws.Range("S3").Value = 0.00001341292715880065
$ws.Range("T3").Value = 0.00001341292715880065

# Row 4
$ws.Range("G4").Value = 2.543891666666667
$ws.Range("H4").Value = 7.631675
$ws.Range("I4").Value = 0.0487891203504456
$ws.Range("J4").Value = 0.04878912035044559
$ws.Range("O4").Value = 0.002616675800765965
$ws.Range("P4").Value = 0.002616675800765965
$ws.Range("Q4").Value = 0.2732029414694445
$ws.Range("R4").Value = 2.458826473225
$ws.Range("S4").Value = 0.0001276653105616693
$ws.Range("T4").Value = 0.0001276653105616693

# Row 5
$ws.Range("G5").Value = 2.543891666666667
$ws.Range("H5").Value = 7.631675
$ws.Range("I5").Value = 0.0487891203504456
$ws.Range("J5").Value = 0.04878912035044559
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.113565
$ws.Range("N5").Value = 0.340695
$ws.Range("O5").Value = 0.002766990480503436
$ws.Range("P5").Value = 0.002766990480503436
$ws.Range("Q5").Value = 0.288897057125
$ws.Range("R5").Value = 2.600073514125
$ws.Range("S5").Value = 0.0001349990315618194
$ws.Range("T5").Value = 0.0001349990315618194

# Row 6
$ws.Range("I6").Value = 0.0400662233111763
$ws.Range("J6").Value = 0.0400662233111763
$ws.Range("M6").Value = 40.81054266666667
$ws.Range("N6").Value = 122.431628
$ws.Range("O6").Value = 0.9943414173631485
$ws.Range("P6").Value = 0.9943414173631485
$ws.Range("Q6").Value = 85.25628442136667
$ws.Range("R6").Value = 767.3065597923
$ws.Range("S6").Value = 0.03983950527562346
$ws.Range("T6").Value = 0.03983950527562347

# Row 7
$ws.Range("I7").Value = 0.0400662233111763
$ws.Range("J7").Value = 0.0400662233111763
$ws.Range("O7").Value = 0.0002749163555820933
$ws.Range("P7").Value = 0.0002749163555820933
$ws.Range("S7").Value = 0.0000110148600946469
$ws.Range("T7").Value = 0.0000110148600946469

# Row 8
$ws.Range("I8").Value = 0.0400662233111763
$ws.Range("J8").Value = 0.0400662233111763
$ws.Range("O8").Value = 0.002616675800765965
$ws.Range("P8").Value = 0.002616675800765965
$ws.Range("S8").Value = 0.0001048403169664402
$ws.Range("T8").Value = 0.0001048403169664402

# Row 9
$ws.Range("I9").Value = 0.0400662233111763
$ws.Range("J9").Value = 0.0400662233111763
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.113565
$ws.Range("N9").Value = 0.340695
$ws.Range("O9").Value = 0.002766990480503436
$ws.Range("P9").Value = 0.002766990480503436
$ws.Range("Q9").Value = 0.237245802375
$ws.Range("R9").Value = 2.135212221375
$ws.Range("S9").Value = 0.0001108628584917497
$ws.Range("T9").Value = 0.0001108628584917497

# Row 10
$ws.Range("G10").Value = 1.800112666666666
$ws.Range("H10").Value = 5.400338
$ws.Range("I10").Value = 0.03452423493074386
$ws.Range("J10").Value = 0.03452423493074386
$ws.Range("M10").Value = 40.81054266666667
$ws.Range("N10").Value = 122.431628
$ws.Range("O10").Value = 0.9943414173631485
$ws.Range("P10").Value = 0.9943414173631485
$ws.Range("Q10").Value = 73.4635747878071
$ws.Range("R10").Value = 661.172173090264
$ws.Range("S10").Value = 0.03432887669441417
$ws.Range("T10").Value = 0.03432887669441417

# Row 11
$ws.Range("G11").Value = 1.800112666666666
$ws.Range("H11").Value = 5.400338
$ws.Range("I11").Value = 0.03452423493074386
$ws.Range("J11").Value = 0.03452423493074386
$ws.Range("O11").Value = 0.0002749163555820933
$ws.Range("P11").Value = 0.0002749163555820933
$ws.Range("Q11").Value = 0.02031127125555555
$ws.Range("R11").Value = 0.1828014413
$ws.Range("S11").Value = 0.000009491276846420106
$ws.Range("T11").Value = 0.000009491276846420108

# Row 12
$ws.Range("G12").Value = 1.800112666666666
$ws.Range("H12").Value = 5.400338
$ws.Range("I12").Value = 0.03452423493074386
$ws.Range("J12").Value = 0.03452423493074386
$ws.Range("O12").Value = 0.002616675800765965
$ws.Range("P12").Value = 0.002616675800765965
$ws.Range("Q12").Value = 0.1933242999117778
$ws.Range("R12").Value = 1.739918699206
$ws.Range("S12").Value = 0.00009033873008323649
$ws.Range("T12").Value = 0.0000903387300832365

# Row 13
$ws.Range("G13").Value = 1.800112666666666
$ws.Range("H13").Value = 5.400338
$ws.Range("I13").Value = 0.03452423493074386
$ws.Range("J13").Value = 0.03452423493074386
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.113565
$ws.Range("N13").Value = 0.340695
$ws.Range("O13").Value = 0.002766990480503436
$ws.Range("P13").Value = 0.002766990480503436
$ws.Range("Q13").Value = 0.20442979499
$ws.Range("R13").Value = 1.83986815491
$ws.Range("S13").Value = 0.00009552822940003247
$ws.Range("T13").Value = 0.00009552822940003247

# Row 14
$ws.Range("G14").Value = 45.70747266666666
$ws.Range("H14").Value = 137.122418
$ws.Range("I14").Value = 0.8766204214076343
$ws.Range("J14").Value = 0.8766204214076342
$ws.Range("M14").Value = 40.81054266666667
$ws.Range("N14").Value = 122.431628
$ws.Range("O14").Value = 0.9943414173631485
$ws.Range("P14").Value = 0.9943414173631485
$ws.Range("Q14").Value = 1865.3467634485
$ws.Range("R14").Value = 16788.1208710365
$ws.Range("S14").Value = 0.8716599923119476
$ws.Range("T14").Value = 0.8716599923119475

# Row 15
$ws.Range("G15").Value = 45.70747266666666
$ws.Range("H15").Value = 137.122418
$ws.Range("I15").Value = 0.8766204214076343
$ws.Range("J15").Value = 0.8766204214076342
$ws.Range("O15").Value = 0.0002749163555820933
$ws.Range("P15").Value = 0.0002749163555820933
$ws.Range("Q15").Value = 0.5157326499222221
$ws.Range("R15").Value = 4.6415938493
$ws.Range("S15").Value = 0.0002409972914822257
$ws.Range("T15").Value = 0.0002409972914822257

# Row 16
$ws.Range("G16").Value = 45.70747266666666
$ws.Range("H16").Value = 137.122418
$ws.Range("I16").Value = 0.8766204214076343
$ws.Range("J16").Value = 0.8766204214076342
$ws.Range("O16").Value = 0.002616675800765965
$ws.Range("P16").Value = 0.002616675800765965
$ws.Range("Q16").Value = 4.908784498685111
$ws.Range("R16").Value = 44.17906048816599
$ws.Range("S16").Value = 0.002293831443154619
$ws.Range("T16").Value = 0.002293831443154619

# Row 17
$ws.Range("G17").Value = 45.70747266666666
$ws.Range("H17").Value = 137.122418
$ws.Range("I17").Value = 0.8766204214076343
$ws.Range("J17").Value = 0.8766204214076342
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.113565
$ws.Range("N17").Value = 0.340695
$ws.Range("O17").Value = 0.002766990480503436
$ws.Range("P17").Value = 0.002766990480503436
$ws.Range("Q17").Value = 5.19076913339
$ws.Range("R17").Value = 46.71692220051
$ws.Range("S17").Value = 0.002425600361049834
$ws.Range("T17").Value = 0.002425600361049834
